$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# NOTE on ordering: we always write the cell's value/formula FIRST, and only
# afterwards copy over the number-format/font/border styling with
# PasteSpecial(xlPasteFormats). Doing it the other way around (format the
# blank cell, then overwrite its value/formula) leaves pre-existing
# whole-column formulas like D2's =SUM(B:B) stuck on their old cached
# result in this runtime, even though the cell's own value reads back fine.
$xlPasteFormats = -4122

# --- Row 69: new journal entry, date 2024-09-02 (serial 45537) ---
# Formatting mirrors row 25 (a "new year" separator row): A s=4 (bold +
# bottom border date format), B/C s=5 (bold + bottom border).
$ws.Range("A69").Value = 45537
$ws.Range("B69").Formula = "=3+3"
$ws.Range("C69").Formula = "=C68+B69"

$ws.Range("A25").Copy()
$ws.Range("A69").PasteSpecial($xlPasteFormats)
$ws.Range("B25").Copy()
$ws.Range("B69").PasteSpecial($xlPasteFormats)
$ws.Range("C25").Copy()
$ws.Range("C69").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Row 70: date 2024-09-06 (serial 45541) ---
# Plain data-row formatting, mirrors row 26 (A s=1 date, B/C unstyled).
$ws.Range("A70").Value = 45541
$ws.Range("B70").Value = 3
$ws.Range("C70").Formula = "=C69+B70"

$ws.Range("A26").Copy()
$ws.Range("A70").PasteSpecial($xlPasteFormats)
$ws.Range("B26").Copy()
$ws.Range("B70").PasteSpecial($xlPasteFormats)
$ws.Range("C26").Copy()
$ws.Range("C70").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Row 71: date 2024-09-07 (serial 45542) ---
$ws.Range("A71").Value = 45542
$ws.Range("B71").Value = 3
$ws.Range("C71").Formula = "=C70+B71"

$ws.Range("A26").Copy()
$ws.Range("A71").PasteSpecial($xlPasteFormats)
$ws.Range("B26").Copy()
$ws.Range("B71").PasteSpecial($xlPasteFormats)
$ws.Range("C26").Copy()
$ws.Range("C71").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Move selection down to the new bottom of the journal, like the author
# did after entering the new rows.
[void]$ws.Range("B72").Select()

$excel.Calculate()
